$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Random_Forest
$ws.Range("D3").Value = 99.03
$ws.Range("E3").Value = 98.98
$ws.Range("G3").Value = 99.5

# Row 4 - Deep_Neural_Network
$ws.Range("C4").Value = 94.38
$ws.Range("D4").Value = 87.39
$ws.Range("E4").Value = 85.70999999999999
$ws.Range("F4").Value = 95.09999999999999
$ws.Range("G4").Value = 90.5
